$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows where Target cluster = "ECs" (old rows 2, 4, 6), keeping only
# rows where Target cluster = "MuSCs" (old rows 3, 5, 7 become new rows 2, 3, 4).
# Delete bottom-to-top so row indices of not-yet-deleted rows stay stable.
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(2).Delete()

# Row 2: ECs -> MuSCs, with recomputed TPM-derived values
$ws.Range("D2").Value = "MuSCs"
$ws.Range("G2").Value = 77.127561
$ws.Range("H2").Value = 231.382683
$ws.Range("I2").Value = 0.2899056040435161
$ws.Range("J2").Value = 0.2899056040435161
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1935556666666667
$ws.Range("N2").Value = 0.5806669999999999
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 14.928476487729
$ws.Range("R2").Value = 134.356288389561
$ws.Range("S2").Value = 0.2899056040435161
$ws.Range("T2").Value = 0.2899056040435161

# Row 3: FAPs -> MuSCs, with recomputed TPM-derived values
$ws.Range("G3").Value = 117.8650183333333
$ws.Range("H3").Value = 353.595055
$ws.Range("I3").Value = 0.443028781054351
$ws.Range("J3").Value = 0.443028781054351
$ws.Range("N3").Value = 0.5806669999999999
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 22.81344220018722
$ws.Range("R3").Value = 205.320979801685
$ws.Range("S3").Value = 0.443028781054351
$ws.Range("T3").Value = 0.443028781054351

# Row 4: MuSCs -> MuSCs, with recomputed TPM-derived values
$ws.Range("G4").Value = 71.05112566666666
$ws.Range("H4").Value = 213.153377
$ws.Range("I4").Value = 0.267065614902133
$ws.Range("J4").Value = 0.2670656149021329
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.1935556666666667
$ws.Range("N4").Value = 0.5806669999999999
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 13.75234799582878
$ws.Range("R4").Value = 123.771131962459
$ws.Range("S4").Value = 0.267065614902133
$ws.Range("T4").Value = 0.2670656149021329
